$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "56.949.66"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +2.52%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.008.89"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +1.62%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E4").Value = "  -0.27%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "513.48"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +4.50%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "139.87"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +5.53%  "

$ws.Range("E7").Value = "  -0.08%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.436"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +4.06%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "7.48"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +5.00%  "

$ws.Range("E10").Value = "  +7.49%  "

$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("E12").Value = "  +1.77%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "3.519.58"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +1.24%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "26.01"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +5.02%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0000157"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +11.96%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "56.890.75"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.46%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.000.71"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +1.28%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "5.97"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +5.12%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.59"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +3.71%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.87"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +5.05%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "328.14"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +2.88%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("E23").Value = "  +5.38%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "63.48"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +5.25%  "

$ws.Range("E25").Value = "  +5.09%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.10%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0918"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +8.79%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "6.70"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.93%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.08"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +8.22%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.24"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +5.89%  "

$ws.Range("E31").Value = "  +6.66%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "20.55"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +5.58%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "155.47"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +3.84%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.59"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +4.45%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.73"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.87%  "

$ws.Range("E36").Value = "  -1.40%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0681"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +4.88%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "23.95"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +2.57%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.037.08"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.35%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "37.17"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("E41").Value = "  -0.23%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.298.01"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +7.76%  "

$ws.Range("E43").Value = "  +2.71%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.70"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +4.65%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("E46").Value = "  +3.55%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.97"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +10.54%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "5.90"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +5.93%  "

$ws.Range("E49").Value = "  +1.79%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "19.48"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("E51").Value = "  +4.72%  "
